$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string (A1)
$ws.Range("A1").Value = "Datos actualizados a 23 de Agosto de 2020 a las 02:49"

# Row 4: Estados Unidos - refreshed totals
$ws.Range("B4").Value = 5840869
$ws.Range("C4").Value = 43270
$ws.Range("D4").Value = 3148080
$ws.Range("E4").Value = 2512621
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 968
$ws.Range("H4").Value = 180168

# Row 5: Brasil - refreshed totals
$ws.Range("B5").Value = 3582698
$ws.Range("C5").Value = 46210
$ws.Range("D5").Value = 2709638
$ws.Range("E5").Value = 758783
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 823
$ws.Range("H5").Value = 114277

# Rows 113-114: Montenegro overtakes Nicaragua in ranking order
$ws.Range("A113").Value = "Montenegro"
$ws.Range("B113").Value = 4313
$ws.Range("C113").Value = 36
$ws.Range("D113").Value = 3333
$ws.Range("E113").Value = 896
$ws.Range("F113").Value = 0
$ws.Range("G113").Value = 2
$ws.Range("H113").Value = 84

$ws.Range("A114").Value = "Nicaragua"
$ws.Range("B114").Value = 4311
$ws.Range("C114").Value = 0
$ws.Range("D114").Value = 2913
$ws.Range("E114").Value = 1265
$ws.Range("F114").Value = 0
$ws.Range("G114").Value = 0
$ws.Range("H114").Value = 133

# Rows 175-176: San Martin (Parte Holandesa) overtakes Papua Nueva Guinea
$ws.Range("A175").Value = "San Martin (Parte Holandesa)"
$ws.Range("B175").Value = 368
$ws.Range("C175").Value = 15
$ws.Range("D175").Value = 147
$ws.Range("E175").Value = 204
$ws.Range("F175").Value = 0
$ws.Range("G175").Value = 0
$ws.Range("H175").Value = 17

$ws.Range("A176").Value = "Papua Nueva Guinea"
$ws.Range("B176").Value = 361
$ws.Range("C176").Value = 0
$ws.Range("D176").Value = 198
$ws.Range("E176").Value = 159
$ws.Range("F176").Value = 0
$ws.Range("G176").Value = 0
$ws.Range("H176").Value = 4

# Row 188: Bermudas - refreshed totals
$ws.Range("D188").Value = 149
$ws.Range("E188").Value = 9

# Row 189: Barbados - refreshed totals
$ws.Range("B189").Value = 158
$ws.Range("C189").Value = 1
$ws.Range("D189").Value = 126
$ws.Range("E189").Value = 25
